# Updated cryptos list values (coin name, link, price, 1h volume %) per latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "62.114.73"
    "E2" = "  -0.81%  "
    "D3" = "3.424.14"
    "E3" = "  -0.98%  "
    "E4" = "  -0.16%  "
    "D5" = "408.26"
    "E5" = "  -0.95%  "
    "D6" = "134.18"
    "E6" = "  +4.52%  "
    "E7" = "  +0.03%  "
    "D8" = "1.00"
    "E8" = "  -0.12%  "
    "E9" = "  -1.12%  "
    "E10" = "  -2.81%  "
    "D11" = "42.67"
    "E11" = "  -2.16%  "
    "E12" = "  -1.41%  "
    "B13" = "Polkadot"
    "C13" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D13" = "8.44"
    "E13" = "  -3.04%  "
    "B14" = "Chainlink"
    "C14" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D14" = "19.90"
    "E14" = "  -1.23%  "
    "B15" = "WrappedEther"
    "C15" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D15" = "3.413.11"
    "E15" = "  -0.27%  "
    "B16" = "WrappedBTC"
    "C16" = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
    "D16" = "62.155.09"
    "E16" = "  -0.73%  "
    "B17" = "Polygon"
    "C17" = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
    "D17" = "1.01"
    "E17" = "  -2.91%  "
    "B18" = "Uniswap"
    "C18" = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
    "D18" = "11.00"
    "E18" = "  +0.50%  "
    "B19" = "ShibaInu"
    "C19" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
    "D19" = "0.0000131"
    "E19" = "  -2.58%  "
    "B20" = "ImmutableX"
    "C20" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D20" = "3.21"
    "E20" = "  -4.47%  "
    "B21" = "Litecoin"
    "C21" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "D21" = "84.72"
    "E21" = "  +3.83%  "
    "B22" = "BitcoinCash"
    "C22" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "D22" = "315.25"
    "E22" = "  +1.32%  "
    "B23" = "InternetComputer(DFINITY)"
    "C23" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D23" = "12.84"
    "E23" = "  -2.64%  "
    "B24" = "PancakeSwap"
    "C24" = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
    "D24" = "3.19"
    "E24" = "  +0.06%  "
    "B25" = "LEO"
    "C25" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D25" = "4.78"
    "E25" = "  +9.45%  "
    "B26" = "EthereumClassic"
    "C26" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D26" = "29.71"
    "E26" = "  -1.70%  "
    "B27" = "Filecoin"
    "C27" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D27" = "8.21"
    "E27" = "  +2.02%  "
    "B28" = "RenderToken"
    "C28" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D28" = "7.67"
    "E28" = "  -1.80%  "
    "B29" = "Toncoin"
    "C29" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "D29" = "2.75"
    "E29" = "  +4.55%  "
    "B30" = "Kaspa"
    "C30" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D30" = "0.175"
    "E30" = "  -1.27%  "
    "B31" = "Hedera"
    "C31" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D31" = "0.116"
    "E31" = "  -4.32%  "
    "B32" = "InjectiveProtocol"
    "C32" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D32" = "42.82"
    "E32" = "  -4.24%  "
    "B33" = "Dai"
    "C33" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D33" = "1.00"
    "E33" = "  -0.04%  "
    "B34" = "Cosmos"
    "C34" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D34" = "11.40"
    "E34" = "  -6.37%  "
    "B35" = "VeChain"
    "C35" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D35" = "0.0485"
    "E35" = "  -2.30%  "
    "B36" = "OKB"
    "C36" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D36" = "52.09"
    "E36" = "  -1.15%  "
    "B37" = "FirstDigitalUSD"
    "C37" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D37" = "1.00"
    "E37" = "  +0.17%  "
    "B38" = "LidoDAOToken"
    "C38" = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
    "D38" = "3.43"
    "E38" = "  -3.74%  "
    "B39" = "Stacks"
    "C39" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D39" = "2.98"
    "E39" = "  -1.15%  "
    "B40" = "ARBITRUM"
    "C40" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D40" = "2.00"
    "E40" = "  -0.20%  "
    "B41" = "Monero"
    "C41" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D41" = "137.67"
    "E41" = "  -0.05%  "
    "B42" = "Stellar"
    "C42" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D42" = "0.125"
    "E42" = "  +0.15%  "
    "B43" = "TheGraph"
    "C43" = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
    "D43" = "0.297"
    "E43" = "  +3.53%  "
    "B44" = "NEARProtocol"
    "C44" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D44" = "4.02"
    "E44" = "  +0.88%  "
    "B45" = "Celestia"
    "C45" = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
    "D45" = "16.84"
    "E45" = "  -5.67%  "
    "B46" = "WEMIXToken"
    "C46" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D46" = "2.22"
    "E46" = "  -2.18%  "
    "B47" = "EnergySwap"
    "C47" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D47" = "21.47"
    "E47" = "  -4.31%  "
    "B48" = "Maker"
    "C48" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D48" = "2.130.99"
    "E48" = "  -4.33%  "
    "B49" = "ApeXProtocol"
    "C49" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "D49" = "2.28"
    "E49" = "  -4.69%  "
    "B50" = "ThetaToken"
    "C50" = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    "D50" = "1.93"
    "E50" = "  +3.39%  "
    "B51" = "BEAM"
    "C51" = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
    "D51" = "0.0352"
    "E51" = "  +4.05%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text formatting so numeric-looking strings (e.g. "1.00", "408.26") are not
    # coerced into actual numbers and keep their original text representation.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
